$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (MAE)
$ws.Range("B2").Value = 1.262
$ws.Range("C2").Value = 1.141
$ws.Range("D2").Value = 0.791
$ws.Range("E2").Value = 0.6909999999999999
$ws.Range("F2").Value = 1.509

# Row 3 (MSE)
$ws.Range("B3").Value = 2.507
$ws.Range("C3").Value = 2.335
$ws.Range("D3").Value = 1.137
$ws.Range("E3").Value = 0.833
$ws.Range("F3").Value = 4.191

# Row 4 (R2)
$ws.Range("B4").Value = 0.659
$ws.Range("C4").Value = 0.8
$ws.Range("D4").Value = 0.789
$ws.Range("E4").Value = 0.645
$ws.Range("F4").Value = 0.856
